$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 23 values (W(kg) and Numb columns)
$ws.Range("G23").Value = 0.304
$ws.Range("H23").Value = 39

# Delete row 37 entirely (the "Parapenaeus longirostris" / PAPELON duplicate
# entry) - this shifts rows 38-40 up to become 37-39, so the previously
# last row (40, "Wood NA") disappears and the sheet shrinks to 39 rows.
$ws.Rows("37").Delete()
